$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Datos")
$ws.Activate()

# Update the e-mail value in N2 (was jruav@devco.com.co -> jrua@todo1.net),
# keep existing hyperlink + style (s="9") untouched.
$ws.Range("N2").Value = "jrua@todo1.net"

# Add the new "tipoCorreo" column header/value (O1/O2).
$ws.Range("O1").Value = "tipoCorreo"
$ws.Range("O2").Value = "Laboral"
$ws.Range("O2").NumberFormat = "@"

# Rename the N column header from "correoUsuario" to "correo".
$ws.Range("N1").Value = "correo"

# Add the new "numeroCelular" column (P1/P2).
$ws.Range("P1").Value = "numeroCelular"
$ws.Range("P2").Value = 3146834995

# Match the new selection / view state.
$ws.Range("P2").Select()
